$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2"  = 0.08707909936379449
    "C2"  = 0.1849415347999282
    "D2"  = 0.6087619768527019
    "B3"  = 0.2343920533817323
    "C3"  = 0.5621024059030257
    "B4"  = 0.7136637901383804
    "B5"  = 0.6611739089257214
    "C5"  = 0.2424661772861064
    "D5"  = 0.343437461374561
    "E5"  = 0.2316648163334098
    "B6"  = 0.52173245889808
    "C6"  = 0.3916150992080441
    "D6"  = 0.0731403976960974
    "B7"  = 0.6205335770846683
    "C7"  = 0.1160860902143776
    "B8"  = 0.3924296930042493
    "B9"  = 0.4747974866332783
    "C9"  = -0.06986322861814367
    "D9"  = 0.2092250849631835
    "E9"  = 0.02501891626858249
    "B10" = 0.2191407266580338
    "C10" = 0.2222972325179511
    "D10" = -0.09185372756793091
    "B11" = 0.224871884804233
    "C11" = -0.0167560416292489
    "B12" = 0.1712949347374914
    "B13" = -0.01257980384821508
    "C13" = 0.1395187281562265
    "D13" = 0.1036987463175434
    "B14" = 0.2732997081603292
    "C14" = 0.04948772950879379
    "B15" = 0.1075450023515467
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
